# semana 29 de 2025
# Adds a new "week 29" column (AF) to the weekly IRA-hospital revision sheet,
# mirroring the existing week columns (D..AE hold weeks 1..28), and fills in
# the week-29 counts per UPGD row. Also adds the missing UPGD name for the
# newly-reporting institution in row 54.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new week column (AF1 = "29"), matching the style of the
# other header cells in row 1 (bold / centered -> same style as A1).
# Copy the header style first, then force text storage (NumberFormat "@")
# before assigning, otherwise the numeric-looking string "29" gets
# auto-coerced to a number.
$ws.Cells.Item(1, 32).Style = $ws.Cells.Item(1, 31).Style
$ws.Cells.Item(1, 32).NumberFormat = "@"
$ws.Cells.Item(1, 32).Value = "29"

# New UPGD name that showed up for the first time this week.
$ws.Cells.Item(54, 3).Value = "CLINICA MEDICA TURIN SAS"

# Week 29 (column AF = column 32) values per row.
$af = @{
  2  = 0
  3  = 0
  4  = 0
  5  = 0
  6  = 34
  7  = 1
  8  = 28
  9  = 0
  12 = 0
  14 = 0
  15 = 0
  17 = 0
  23 = 0
  24 = 0
  25 = 1
  26 = 0
  27 = 0
  28 = 69
  29 = 4
  30 = 19
  31 = 0
  32 = 0
  34 = 0
  35 = 13
  36 = 0
  37 = 0
  38 = 0
  39 = 0
  40 = 0
  41 = 0
  42 = 0
  43 = 0
  44 = 0
  45 = 0
  46 = 0
  47 = 0
  48 = 0
  49 = 0
  50 = 0
  52 = 0
  53 = 0
  54 = 0
  55 = 0
  56 = 0
  57 = 0
}

foreach ($row in $af.Keys) {
  $ws.Cells.Item($row, 32).Value = $af[$row]
}
